# Applies the cryptos.xlsx price/volume/coin refresh described by the commit.
# Column D ("Price") values are plain digit/decimal text that Excel would
# otherwise auto-coerce into numbers (dropping trailing zeros, flipping tiny
# values to scientific notation, etc). Set-Text forces the literal string via
# a temporary Text number format, then clears that transient formatting so no
# stray style is left behind on cells that were previously unstyled.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Text($range, $text) {
    if ($text -match "^-?\d+(\.\d+)?$") {
        $range.NumberFormat = "@"
        $range.Value = $text
        $range.ClearFormats()
    } else {
        $range.Value = $text
    }
}

Set-Text $ws.Range("D2") "66.707.65"
Set-Text $ws.Range("E2") "  +2.76%  "
Set-Text $ws.Range("D3") "3.688.19"
Set-Text $ws.Range("E3") "  +4.95%  "
Set-Text $ws.Range("D4") "0.999"
Set-Text $ws.Range("E4") "  +0.09%  "
Set-Text $ws.Range("D5") "418.57"
Set-Text $ws.Range("E5") "  -0.24%  "
Set-Text $ws.Range("D6") "130.20"
Set-Text $ws.Range("E6") "  -1.28%  "
Set-Text $ws.Range("D7") "3.681.84"
Set-Text $ws.Range("E7") "  +4.91%  "
Set-Text $ws.Range("D8") "0.641"
Set-Text $ws.Range("E8") "  -0.25%  "
Set-Text $ws.Range("E9") "  -0.04%  "
Set-Text $ws.Range("D10") "0.759"
Set-Text $ws.Range("E10") "  -3.93%  "
Set-Text $ws.Range("D11") "0.181"
Set-Text $ws.Range("E11") "  +8.80%  "
Set-Text $ws.Range("D12") "0.0000389"
Set-Text $ws.Range("E12") "  +45.51%  "
Set-Text $ws.Range("D13") "42.78"
Set-Text $ws.Range("E13") "  -1.56%  "
Set-Text $ws.Range("D14") "10.59"
Set-Text $ws.Range("E14") "  +6.88%  "
Set-Text $ws.Range("D15") "4.258.59"
Set-Text $ws.Range("E15") "  +4.98%  "
Set-Text $ws.Range("E16") "  -0.79%  "
Set-Text $ws.Range("D17") "3.803.75"
Set-Text $ws.Range("E17") "  +7.89%  "
Set-Text $ws.Range("D18") "20.55"
Set-Text $ws.Range("E18") "  -0.18%  "
Set-Text $ws.Range("D19") "12.90"
Set-Text $ws.Range("E19") "  +3.70%  "
Set-Text $ws.Range("E20") "  +1.30%  "
Set-Text $ws.Range("D21") "66.658.27"
Set-Text $ws.Range("E21") "  +3.02%  "
Set-Text $ws.Range("D22") "440.31"
Set-Text $ws.Range("E22") "  -4.78%  "
Set-Text $ws.Range("D23") "16.66"
Set-Text $ws.Range("E23") "  +23.16%  "
Set-Text $ws.Range("D24") "89.21"
Set-Text $ws.Range("E24") "  -2.39%  "
Set-Text $ws.Range("E25") "  -5.07%  "
Set-Text $ws.Range("E26") "  +8.50%  "
Set-Text $ws.Range("D27") "10.31"
Set-Text $ws.Range("E27") "  +0.89%  "
Set-Text $ws.Range("E28") "  -2.01%  "
Set-Text $ws.Range("E29") "  +3.91%  "
Set-Text $ws.Range("B30") "Hedera"
Set-Text $ws.Range("C30") "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-Text $ws.Range("D30") "0.124"
Set-Text $ws.Range("E30") "  +8.95%  "
Set-Text $ws.Range("B31") "Cosmos"
Set-Text $ws.Range("C31") "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-Text $ws.Range("D31") "12.70"
Set-Text $ws.Range("E31") "  +1.03%  "
Set-Text $ws.Range("E32") "  +2.36%  "
Set-Text $ws.Range("E33") "  -5.37%  "
Set-Text $ws.Range("D34") "0.165"
Set-Text $ws.Range("E34") "  -1.25%  "
Set-Text $ws.Range("D35") "41.33"
Set-Text $ws.Range("E35") "  +1.75%  "
Set-Text $ws.Range("D36") "57.03"
Set-Text $ws.Range("E36") "  -1.79%  "
Set-Text $ws.Range("E37") "  -0.01%  "
Set-Text $ws.Range("D38") "0.0492"
Set-Text $ws.Range("E38") "  -3.86%  "
Set-Text $ws.Range("B39") "ThetaToken"
Set-Text $ws.Range("C39") "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-Text $ws.Range("D39") "3.17"
Set-Text $ws.Range("E39") "  +35.69%  "
Set-Text $ws.Range("B40") "PEPE"
Set-Text $ws.Range("C40") "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-Text $ws.Range("D40") "0.0₃0744"
Set-Text $ws.Range("E40") "  +6.97%  "
Set-Text $ws.Range("B41") "Stellar"
Set-Text $ws.Range("C41") "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-Text $ws.Range("D41") "0.147"
Set-Text $ws.Range("E41") "  +2.49%  "
Set-Text $ws.Range("B42") "EnergySwap"
Set-Text $ws.Range("C42") "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-Text $ws.Range("D42") "28.60"
Set-Text $ws.Range("E42") "  +30.31%  "
Set-Text $ws.Range("D43") "0.997"
Set-Text $ws.Range("E43") "  -0.10%  "
Set-Text $ws.Range("E44") "  +0.57%  "
Set-Text $ws.Range("D45") "148.74"
Set-Text $ws.Range("E45") "  +2.06%  "
Set-Text $ws.Range("B46") "ARBITRUM"
Set-Text $ws.Range("C46") "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-Text $ws.Range("D46") "2.10"
Set-Text $ws.Range("E46") "  +2.91%  "
Set-Text $ws.Range("B47") "NEARProtocol"
Set-Text $ws.Range("C47") "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-Text $ws.Range("D47") "4.40"
Set-Text $ws.Range("E47") "  -3.33%  "
Set-Text $ws.Range("D48") "2.89"
Set-Text $ws.Range("E48") "  -8.27%  "
Set-Text $ws.Range("D49") "2.58"
Set-Text $ws.Range("E49") "  -7.14%  "
Set-Text $ws.Range("D50") "0.305"
Set-Text $ws.Range("E50") "  -5.35%  "
Set-Text $ws.Range("E51") "  +11.84%  "
